# Atualizações após dados finais recebidos
# O professor enviou os últimos dados em 24/03/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Fix existing values (C12 and C13) ---
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(13, 3).Value = 2

# --- Append new rows 32-39 (MS-31 .. MS-38) ---
$newData = @(
    @(31, "MS-31", 2),
    @(32, "MS-32", 2),
    @(33, "MS-33", 1),
    @(34, "MS-34", 1),
    @(35, "MS-35", 2),
    @(36, "MS-36", 2),
    @(37, "MS-37", 1),
    @(38, "MS-38", 2)
)

$row = 32
foreach ($entry in $newData) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# --- Update selection to reflect new data range ---
$ws.Range("C2:C39").Select()
